# Fruta / hortaliza, semanal
# Insert two new weekly price rows at the top of the data block (rows 28-29),
# pushing the existing rows 28-69 down to rows 30-71.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 28..69 down by two rows (new rows 28:29 are blank).
$ws.Rows("28:29").Insert()

# New row 28: Clementina / Primera, Region de O'Higgins
$ws.Cells.Item(28, 1).Value = 7
$ws.Cells.Item(28, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(28, 3).Value = "Ñuble"
$ws.Cells.Item(28, 4).Value = 44413
$ws.Cells.Item(28, 5).Value = 16
$ws.Cells.Item(28, 6).Value = "Fruta"
$ws.Cells.Item(28, 7).Value = 100102
$ws.Cells.Item(28, 8).Value = "Cítricos"
$ws.Cells.Item(28, 9).Value = 100102004
$ws.Cells.Item(28, 10).Value = "Mandarina"
$ws.Cells.Item(28, 11).Value = "Clementina"
$ws.Cells.Item(28, 12).Value = "Primera"
$ws.Cells.Item(28, 13).Value = 100
$ws.Cells.Item(28, 14).Value = 6000
$ws.Cells.Item(28, 15).Value = 6500
$ws.Cells.Item(28, 16).Value = 6250
$ws.Cells.Item(28, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(28, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(28, 19).Value = 625
$ws.Cells.Item(28, 20).Value = 10

# New row 29: Clementina / Segunda, Region de O'Higgins
$ws.Cells.Item(29, 1).Value = 7
$ws.Cells.Item(29, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(29, 3).Value = "Ñuble"
$ws.Cells.Item(29, 4).Value = 44413
$ws.Cells.Item(29, 5).Value = 16
$ws.Cells.Item(29, 6).Value = "Fruta"
$ws.Cells.Item(29, 7).Value = 100102
$ws.Cells.Item(29, 8).Value = "Cítricos"
$ws.Cells.Item(29, 9).Value = 100102004
$ws.Cells.Item(29, 10).Value = "Mandarina"
$ws.Cells.Item(29, 11).Value = "Clementina"
$ws.Cells.Item(29, 12).Value = "Segunda"
$ws.Cells.Item(29, 13).Value = 100
$ws.Cells.Item(29, 14).Value = 5000
$ws.Cells.Item(29, 15).Value = 5500
$ws.Cells.Item(29, 16).Value = 5250
$ws.Cells.Item(29, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(29, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(29, 19).Value = 525
$ws.Cells.Item(29, 20).Value = 10
